$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Institute Address Line 1 changes, Institute Address Line 2 is cleared,
# Member Country changes from USA to Spain.
$ws.Range("K4").Value = "abc street"
$ws.Range("L4").Value = ""
$ws.Range("P4").Value = "Spain"

# Update the active selection on the sheet (matches the new cursor position
# recorded in the saved view state).
$ws.Range("P5").Select()
